$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: bump the count in F1 (was 2 data rows tracked, now 6 total incl. header context)
$ws.Range("F1").Value = 6

# Row 2 (already existed) - keep as-is but rewrite explicitly for consistency
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"
$ws.Range("C2").Value = 0.971200102199074
$ws.Range("D2").Value = "сандальки"
$ws.Range("E2").Value = 300

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"
$ws.Range("C3").Value = 0.9765649613078703
$ws.Range("D3").Value = "мозгииии"
$ws.Range("E3").Value = 1000

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"
$ws.Range("C4").Value = 0.9845908482523147
$ws.Range("D4").Value = "илюха бесценен"
$ws.Range("E4").Value = 20342390

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"
$ws.Range("C5").Value = 0.9938976919791667
$ws.Range("D5").Value = "а"
$ws.Range("E5").Value = 1

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "9AE7AE27BD03B73BBF7E1D495899D72BADEB2BFC476959783685DD330628990D"
$ws.Range("C6").Value = 0.9946925317708334
$ws.Range("D6").Value = "бээ"
$ws.Range("E6").Value = 1

# Apply the time format (h:mm:ss, numFmtId 21) to the Datetime column rows 2-6
$ws.Range("C2:C6").NumberFormat = "h:mm:ss"
